{"js": "// Add the missing surname/given-name (\"Mantura\" / \"Bertilo\") to the\n// \"Cueva ... Alvaro\" responsible-person entries so the full name reads\n// \"Cueva Mantura, Alvaro Bertilo\" (and, in the shorter revision-history\n// table, \"Cueva Mantura Alvaro Bertilo\").\n\nconst body = context.document.body;\n\n// 1) \"RESPONSABLES\" table: a standalone paragraph whose whole text is\n//    \"Cueva\" becomes \"Cueva Mantura, Alvaro Bertilo\".\nconst firstResults = body.search(\"Cueva\", { matchCase: true, matchWholeWord: true });\nfirstResults.load(\"items\");\nawait context.sync();\n\nlet firstTarget = null;\nfor (const r of firstResults.items) {\n  r.load(\"text\");\n}\nawait context.sync();\nfor (const r of firstResults.items) {\n  if (r.text === \"Cueva\") {\n    firstTarget = r;\n    break;\n  }\n}\nif (firstTarget) {\n  firstTarget.insertText(\"Cueva Mantura, Alvaro Bertilo\", \"Replace\");\n  await context.sync();\n}\n\n// 2) \"HISTORIAL DE REVISIONES\" table: the run sequence reading\n//    \"Cueva Alvaro\" becomes \"Cueva Mantura Alvaro Bertilo\".\nconst secondResults = body.search(\"Cueva Alvaro\", { matchCase: true });\nsecondResults.load(\"items\");\nawait context.sync();\n\nif (secondResults.items.length > 0) {\n  secondResults.items[0].insertText(\"Cueva Mantura Alvaro Bertilo\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Add the missing surname/given-name (\"Mantura\" / \"Bertilo\") to the\n# \"Cueva ... Alvaro\" responsible-person entries so the full name reads\n# \"Cueva Mantura, Alvaro Bertilo\" (and, in the shorter revision-history\n# table, \"Cueva Mantura Alvaro Bertilo\").\n\n$d = $word.ActiveDocument\n\n# 1) \"RESPONSABLES\" table: a standalone paragraph whose whole text is\n#    \"Cueva\" becomes \"Cueva Mantura, Alvaro Bertilo\".\n$r1 = $d.Content\n$r1.Find.ClearFormatting()\n$r1.Find.Replacement.ClearFormatting()\n$r1.Find.Text = \"Cueva\"\n$r1.Find.MatchWholeWord = $true\n$r1.Find.MatchCase = $true\n$r1.Find.Replacement.Text = \"Cueva Mantura, Alvaro Bertilo\"\n$r1.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 1)\n\n# 2) \"HISTORIAL DE REVISIONES\" table: the run sequence reading\n#    \"Cueva Alvaro\" becomes \"Cueva Mantura Alvaro Bertilo\".\n$r2 = $d.Content\n$r2.Find.ClearFormatting()\n$r2.Find.Replacement.ClearFormatting()\n$r2.Find.Text = \"Cueva Alvaro\"\n$r2.Find.MatchCase = $true\n$r2.Find.Replacement.Text = \"Cueva Mantura Alvaro Bertilo\"\n$r2.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 1)\n"}
